$wb = $excel.ActiveWorkbook

# --- 1. Remove the empty B2/B3 cells on the "ODI Batting" sheet ---
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBatting.Cells.Item(2, 2).ClearContents()
$odiBatting.Cells.Item(3, 2).ClearContents()

# --- 2. Add a new worksheet "ODI Batting Extra" after the last sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "ODI Batting Extra"

# Helper: write a value that must stay a text cell even if it looks numeric
function Set-TextValue($sheet, $row, $col, $val) {
    $cell = $sheet.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

# Header row (row 1)
$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Reuse the same bold/bordered header style already used on the other sheets
$headerStyleSource = $wb.Worksheets.Item("Player Info").Range("A1:D1")
$headerStyleSource.Copy()
$newSheet.Range("A1:F1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Row 2 : 4608 | 10 | | | | NO
Set-TextValue $newSheet 2 1 "4608"
$newSheet.Cells.Item(2, 2).Value = 10
$newSheet.Cells.Item(2, 6).Value = "NO"

# Row 3 : 4625 | 10 | | | | NO
Set-TextValue $newSheet 3 1 "4625"
$newSheet.Cells.Item(3, 2).Value = 10
$newSheet.Cells.Item(3, 6).Value = "NO"

# Row 4 : 4697 | 10 | 0 | 0 | | NO
Set-TextValue $newSheet 4 1 "4697"
$newSheet.Cells.Item(4, 2).Value = 10
Set-TextValue $newSheet 4 3 "0"
Set-TextValue $newSheet 4 4 "0"
$newSheet.Cells.Item(4, 6).Value = "NO"

# Keep the originally-active sheet selected (workbook view is unchanged by this edit)
$wb.Worksheets.Item("Player Info").Activate()

Write-Host "Edit complete"
